$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 13 = task "Actualizar diagrama de paquetes." (CU - Registrar pago alumno area)
# Mark the task status as done and register 1 consumed hour for Día 6 (column W).
$ws.Range("F13").Value = "Hecho"
$ws.Range("W13").Value = 1

# Reflect the new active selection in the view (matches the scrolled/selected state
# left behind after making the edit above).
$ws.Activate()
$ws.Range("W13").Select()
